$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '39.438.85'
$ws.Range("E2").Value = '  +1.56%  '

$ws.Range("D3").Value = '2.158.57'
$ws.Range("E3").Value = '  +3.09%  '

$ws.Range("E4").Value = '  -0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.78'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.39%  '

$ws.Range("E6").Value = '  +1.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.18'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.97%  '

$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.395'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.51%  '

$ws.Range("E10").Value = '  +2.68%  '

$ws.Range("E11").Value = '  -0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.18'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +7.89%  '

$ws.Range("D13").Value = '2.478.07'
$ws.Range("E13").Value = '  +3.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '22.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.87%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.820'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.81%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.57'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.81%  '

$ws.Range("D17").Value = '2.164.97'
$ws.Range("E17").Value = '  +3.38%  '

$ws.Range("D18").Value = '39.436.86'
$ws.Range("E18").Value = '  +1.77%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.44'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.99%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.18'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.17%  '

$ws.Range("E21").Value = '  +1.81%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '229.22'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.80%  '

$ws.Range("E23").Value = '  -0.01%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.48%  '

$ws.Range("E25").Value = '  +1.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.54%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '172.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.16%  '

$ws.Range("E28").Value = '  -0.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.43'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.93%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.65'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.42%  '

$ws.Range("E31").Value = '  +9.59%  '

$ws.Range("E32").Value = '  +1.34%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.67'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.75%  '

$ws.Range("E34").Value = '  +3.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.20'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.93%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0623'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.62%  '

$ws.Range("E37").Value = '  +1.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.59'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.15%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.00'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.34'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.62%  '

$ws.Range("E41").Value = '  +3.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '103.58'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.72%  '

$ws.Range("D43").Value = '1.542.49'
$ws.Range("E43").Value = '  +0.07%  '

$ws.Range("E44").Value = '  +6.00%  '

$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.11'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +7.07%  '

$ws.Range("B46").Value = 'Cronos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0927'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.50%  '

$ws.Range("E47").Value = '  -0.31%  '

$ws.Range("E48").Value = '  +2.23%  '

$ws.Range("E49").Value = '  +1.05%  '

$ws.Range("B50").Value = 'Celestia'
$ws.Range("C50").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.13'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +27.45%  '

$ws.Range("B51").Value = 'RocketPoolETH'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D51").Value = '2.362.03'
$ws.Range("E51").Value = '  +3.12%  '
